$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows 11-17: column A short labels (copy style from existing A column),
# column B long descriptions (copy style from existing E column)

$ws.Range("A11").Value2 = $ws.Range("A2").Value2
$ws.Range("B11").Value2 = $ws.Range("E2").Value2

$ws.Range("A12").Value2 = "Pedal ratio"
$ws.Range("B12").Value2 = $ws.Range("E3").Value2

$ws.Range("A13").Value2 = $ws.Range("A4").Value2
$ws.Range("B13").Value2 = "Master cylider cross section area"

$ws.Range("A14").Value2 = $ws.Range("A5").Value2
$ws.Range("B14").Value2 = $ws.Range("E5").Value2

$ws.Range("A15").Value2 = $ws.Range("A6").Value2
$ws.Range("B15").Value2 = $ws.Range("E6").Value2

$ws.Range("A16").Value2 = $ws.Range("A7").Value2
$ws.Range("B16").Value2 = $ws.Range("E7").Value2

$ws.Range("A17").Value2 = $ws.Range("A8").Value2
$ws.Range("B17").Value2 = $ws.Range("E8").Value2

# Apply the same cell styles as the existing table (A column style "Entrada", E/B column style "Texto Explicativo")
$ws.Range("A11:A17").Style = $ws.Range("A2").Style
$ws.Range("B11:B17").Style = $ws.Range("E2").Style

# Update selection to match the new active range
$ws.Range("A11:B17").Select()
